# "Crear cuenta bancaria terminada"
# Insert a new "codigo" row (with two-digit bank codes 01..14 plus a
# trailing numeric id 15) above the existing "Incorrecto" validation row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 3; this pushes the old row 3 ("Incorrecto", 123,
# "banco*") down to row 4 and shifts all its cell references along with it.
$ws.Rows("3:3").Insert()

# New column D gets a custom width in the final sheet.
$ws.Columns("D").ColumnWidth = 12.28

# Format B3:O3 as Text (@) *before* writing the zero-padded numeric-looking
# strings "01".."14" into them, so they are stored as shared-string text
# instead of being coerced into plain numbers.
$ws.Range("B3:O3").NumberFormat = "@"

$ws.Range("A3").Value = "codigo"
$ws.Range("B3").Value = "01"
$ws.Range("C3").Value = "02"
$ws.Range("D3").Value = "03"
$ws.Range("E3").Value = "04"
$ws.Range("F3").Value = "05"
$ws.Range("G3").Value = "06"
$ws.Range("H3").Value = "07"
$ws.Range("I3").Value = "08"
$ws.Range("J3").Value = "09"
$ws.Range("K3").Value = "10"
$ws.Range("L3").Value = "11"
$ws.Range("M3").Value = "12"
$ws.Range("N3").Value = "13"
$ws.Range("O3").Value = "14"

# P3 stays a plain (unstyled) number.
$ws.Range("P3").Value = 15

# Move the view down so row 4 is the top visible row, and land the
# selection on B11, matching the saved sheet view state.
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B11").Select()

# Workbook window was minimized with an updated window size in the source
# edit; reflect that on the application window as well.
$excel.ActiveWindow.WindowState = -4140
